# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) listed accounting periods in descending
# order (2003 .. 1707). This update re-sorts them in ascending order
# (1707 .. 2003) and refreshes the "Valor Mora" (F) and "Salario Basico"
# (G) figures for every one of the 33 worker/period rows (16-48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period list (was descending 2003 -> 1707).
$periods = @(
    "1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstRow = 16
$lastRow = 48
$newSalarioBasico = 781242

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $i = $row - $firstRow

    # Periodo Mora (text, e.g. "1707")
    $ws.Range("E$row").Value = $periods[$i]

    # Valor Mora: first 14 rows (16-29) drop to 29509, remaining rows (30-48) move to 31249
    if ($row -le 29) {
        $ws.Range("F$row").Value = 29509
    } else {
        $ws.Range("F$row").Value = 31249
    }

    # Salario Basico: updated uniformly for every row
    $ws.Range("G$row").Value = $newSalarioBasico
}
